$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.066391229629517
$ws.Range("B1").Value = 1.422667264938354
$ws.Range("C1").Value = 2.301104307174683
$ws.Range("D1").Value = 4.561849594116211
$ws.Range("E1").Value = 1.900399565696716
